$wb = $excel.ActiveWorkbook

# --- Sheet "B": the D column (loaded weight) gets swapped into column A,
#     the old A values move out to D, and the B column's "D-4" formula
#     collapses down into its already-computed static value (column C was
#     already static). This mirrors pasting column D's values over column A
#     and flattening the B helper formulas to numbers. ---
$wsB = $wb.Worksheets.Item("B")

for ($r = 1; $r -le 9; $r++) {
    $oldA = $wsB.Cells.Item($r, 1).Value2
    $oldB = $wsB.Cells.Item($r, 2).Value2
    $oldD = $wsB.Cells.Item($r, 4).Value2

    $wsB.Cells.Item($r, 1).Value = $oldD
    $wsB.Cells.Item($r, 2).Value = $oldB
    $wsB.Cells.Item($r, 4).Value = $oldA
}

# --- Cosmetic re-entry of the still-identical fill-down formulas on the
#     other matrix sheets so they collapse into shared-formula groups again
#     (same formulas / same values, just how the sheet was last re-saved). ---
$wsA = $wb.Worksheets.Item("A")
$wsA.Range("A1:A9").Formula = "=D1+4"
$wsA.Range("B1:B9").Formula = "=D1+2"

$wsPctA = $wb.Worksheets.Item("%A")
$wsPctA.Range("D1:D9").Formula = "=B1+0.025"

$wsPctB = $wb.Worksheets.Item("%B")
$wsPctB.Range("D1:D9").Formula = "=B1+0.025"

# --- View state: "%B" was the active sheet with C1:C9 selected; now its
#     selection just sits at F10, and "B" becomes the active sheet with the
#     whole of column D selected. ---
$wsPctB.Activate() | Out-Null
$wsPctB.Range("F10").Select() | Out-Null

$wsB.Activate() | Out-Null
$wsB.Range("D1:D1048576").Select() | Out-Null
